$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data[2] = @(0.8189485623975088, 0.2061131566122185, 0.01661023264762918, 0.1162816279270515, 0.5193643288226966, 0.07973214163530429, 0.3402561718116885, 1.682421352132479)
$data[3] = @(0.7180603288547331, 0.1833884566707127, 0.0150604587627825, 0.1113171366859547, 0.5154294819496812, 0.07973214163530429, 0.3024497049345456, 1.683880131782246)
$data[4] = @(0.6559750088289036, 0.1693662821470809, 0.01410483301238941, 0.1083944561521477, 0.5135291143290246, 0.07973214163530429, 0.279281046064753, 1.686508227068003)
$data[5] = @(0.630640865163798, 0.1636350846412995, 0.01371441688199582, 0.1072346831978734, 0.5128838973985168, 0.07973214163530429, 0.2698508697148867, 1.688013158407699)
$data[6] = @(0.6264321403124882, 0.1626824039160226, 0.01364952961485244, 0.1070439813154351, 0.5127845492536878, 0.07973214163530429, 0.2682856754249485, 1.688289212488286)
$data[7] = @(0.6556334792267648, 0.1692890577926676, 0.01409957169603615, 0.1083786889746463, 0.5135198902073981, 0.07973214163530429, 0.27915382187917, 1.686526768222535)
$data[8] = @(0.7841919542180449, 0.198292216557661, 0.01607673009343813, 0.1145436191279501, 0.5179002895122053, 0.07973214163530429, 0.3272112240547926, 1.682563773995298)
$data[9] = @(1.035147002297492, 0.2546071468782145, 0.01992062790958471, 0.1276444601686606, 0.5306039136335343, 0.07973214163530429, 0.4218107545084138, 1.688615801237461)
$data[10] = @(1.218788417262374, 0.2956286025206509, 0.02272324652935964, 0.1379085384689702, 0.54247778529664, 0.07973214163530429, 0.4915458647135438, 1.701600527165283)
$data[11] = @(1.302165886916782, 0.3142115150775453, 0.02399333549242044, 0.1427215359026732, 0.5484382938233523, 0.07973214163530429, 0.5233242407134497, 1.709385711006092)
$data[12] = @(1.333714577942089, 0.3212369096483201, 0.02447356517656374, 0.1445651473780742, 0.5507763043676164, 0.07973214163530429, 0.5353660005513348, 1.712605735224457)
$data[13] = @(1.326921114491995, 0.3197243843482909, 0.02437017185877721, 0.1441671515382623, 0.5502691660462062, 0.07973214163530429, 0.5327722373072419, 1.711900121390556)
$data[14] = @(1.304761915839606, 0.3147897319554431, 0.02403285898760998, 0.1428727873517701, 0.5486290191895336, 0.07973214163530429, 0.5243147636560508, 1.709645163661889)
$data[15] = @(1.291185532924317, 0.3117656005383651, 0.02382614975773123, 0.1420827016407031, 0.5476349318981448, 0.07973214163530429, 0.5191353557944609, 1.708299406551845)
$data[16] = @(1.213336136799455, 0.2944125675467149, 0.02264014333565711, 0.1375969218825475, 0.5420995445361569, 0.07973214163530429, 0.4894701897874683, 1.701129701764216)
$data[17] = @(1.165535675899775, 0.2837468334788866, 0.02191130663226204, 0.134882127026863, 0.5388473054654099, 0.07973214163530429, 0.4712857867037599, 1.697213723266003)
$data[18] = @(1.138026947619039, 0.2776048680875647, 0.02149164558428396, 0.1333341643495629, 0.5370292913402039, 0.07973214163530429, 0.4608318252871584, 1.69513804717019)
$data[19] = @(1.128710398603857, 0.275524057884013, 0.02134947860475478, 0.1328123613503607, 0.5364227594858164, 0.07973214163530429, 0.4572931914514839, 1.694465549761162)
$data[20] = @(1.170625695537638, 0.2848829788800629, 0.02198893969337234, 0.1351697202158846, 0.5391880655414099, 0.07973214163530429, 0.4732210066173792, 1.697612285663354)
$data[21] = @(1.311271285800444, 0.3162394751710451, 0.02413195587920569, 0.1432523998165181, 0.5491085705602785, 0.07973214163530429, 0.526798712166098, 1.710300104299336)
$data[22] = @(1.403047656824413, 0.3366652746546492, 0.0255282988190757, 0.1486576641926263, 0.5560639487228229, 0.07973214163530429, 0.5618613458506303, 1.720178275538245)
$data[23] = @(1.35407843681287, 0.3257699272813852, 0.02478344224559947, 0.1457614227133703, 0.5523083993855238, 0.07973214163530429, 0.5431435036798007, 1.714760389128202)
$data[24] = @(1.168324583561343, 0.2843693589588838, 0.02195384378720178, 0.1350396594651784, 0.5390338467764479, 0.07973214163530429, 0.4723460921623399, 1.697431548524634)
$data[25] = @(0.9673836328364018, 0.239433658964316, 0.01888444193986061, 0.1239897549793767, 0.5267234622929706, 0.07973214163530429, 0.3961792111008577, 1.685487169080318)

$cols = @("B","C","D","E","F","H","M","O")

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "{0}{1}" -f $cols[$i], $r
        $ws.Range($addr).Value = $vals[$i]
    }
}

$wb.Save()